$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# --- Row 2 ---
Set-TextValue $ws.Range("D2") "27.953.43"
$ws.Range("E2").Value = "  -0.57%  "

# --- Row 3 ---
Set-TextValue $ws.Range("D3") "1.814.78"
$ws.Range("E3").Value = "  +2.18%  "

# --- Row 4 ---
Set-TextValue $ws.Range("D4") "1.007"
$ws.Range("E4").Value = "  +0.17%  "

# --- Row 5 ---
Set-TextValue $ws.Range("D5") "337.10"
$ws.Range("E5").Value = "  -0.57%  "

# --- Row 6 ---
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.25%  "

# --- Row 7 ---
Set-TextValue $ws.Range("D7") "0.4225"
$ws.Range("E7").Value = "  +10.66%  "

# --- Row 8 ---
Set-TextValue $ws.Range("D8") "0.3508"
$ws.Range("E8").Value = "  +2.68%  "

# --- Row 9 ---
Set-TextValue $ws.Range("D9") "45.54"
$ws.Range("E9").Value = "  -3.10%  "

# --- Row 10 ---
Set-TextValue $ws.Range("D10") "1.146"
$ws.Range("E10").Value = "  +0.23%  "

# --- Row 11 ---
Set-TextValue $ws.Range("D11") "0.07455"
$ws.Range("E11").Value = "  +1.05%  "

# --- Row 12 ---
Set-TextValue $ws.Range("D12") "22.93"
$ws.Range("E12").Value = "  -1.67%  "

# --- Row 13 ---
Set-TextValue $ws.Range("D13") "1.002"
$ws.Range("E13").Value = "  -0.10%  "

# --- Row 14 ---
Set-TextValue $ws.Range("D14") "6.258"
$ws.Range("E14").Value = "  -2.01%  "

# --- Row 15 ---
Set-TextValue $ws.Range("D15") "7.307"
$ws.Range("E15").Value = "  +0.13%  "

# --- Row 16 ---
Set-TextValue $ws.Range("D16") "1.818.24"
$ws.Range("E16").Value = "  +2.32%  "

# --- Row 17 ---
Set-TextValue $ws.Range("D17") "0.00001088"
$ws.Range("E17").Value = "  +1.07%  "

# --- Row 18 ---
Set-TextValue $ws.Range("D18") "0.06688"
$ws.Range("E18").Value = "  +0.42%  "

# --- Row 19 ---
Set-TextValue $ws.Range("D19") "82.23"
$ws.Range("E19").Value = "  -0.09%  "

# --- Row 20 ---
Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  -0.09%  "

# --- Row 21 ---
Set-TextValue $ws.Range("D21") "17.33"
$ws.Range("E21").Value = "  -0.10%  "

# --- Row 22 ---
Set-TextValue $ws.Range("D22") "6.390"
$ws.Range("E22").Value = "  +0.00%  "

# --- Row 23 ---
Set-TextValue $ws.Range("D23") "28.045.28"
$ws.Range("E23").Value = "  -0.28%  "

# --- Row 24 ---
Set-TextValue $ws.Range("D24") "11.86"
$ws.Range("E24").Value = "  -1.94%  "

# --- Row 25 ---
Set-TextValue $ws.Range("D25") "2.401"
$ws.Range("E25").Value = "  +0.38%  "

# --- Row 26 ---
Set-TextValue $ws.Range("D26") "2.471"
$ws.Range("E26").Value = "  +3.15%  "

# --- Row 27 ---
Set-TextValue $ws.Range("D27") "20.72"
$ws.Range("E27").Value = "  +0.29%  "

# --- Row 28 ---
Set-TextValue $ws.Range("D28") "156.19"
$ws.Range("E28").Value = "  +1.52%  "

# --- Row 29 ---
Set-TextValue $ws.Range("D29") "2.025.99"
$ws.Range("E29").Value = "  +2.42%  "

# --- Row 30 ---
Set-TextValue $ws.Range("D30") "1.305"
$ws.Range("E30").Value = "  -8.35%  "

# --- Row 31 ---
Set-TextValue $ws.Range("D31") "132.61"
$ws.Range("E31").Value = "  -1.65%  "

# --- Row 32 ---
Set-TextValue $ws.Range("D32") "4.083"
$ws.Range("E32").Value = "  +1.43%  "

# --- Row 33 ---
Set-TextValue $ws.Range("D33") "5.986"
$ws.Range("E33").Value = "  -1.20%  "

# --- Row 34 ---
Set-TextValue $ws.Range("D34") "0.09164"
$ws.Range("E34").Value = "  +3.42%  "

# --- Row 35 ---
Set-TextValue $ws.Range("D35") "12.33"
$ws.Range("E35").Value = "  -2.87%  "

# --- Row 36 ---
Set-TextValue $ws.Range("D36") "0.02355"
$ws.Range("E36").Value = "  -2.07%  "

# --- Row 37 ---
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D37") "0.6661"
$ws.Range("E37").Value = "  -2.38%  "

# --- Row 38 ---
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D38") "0.06299"
$ws.Range("E38").Value = "  -0.77%  "

# --- Row 39 ---
Set-TextValue $ws.Range("D39") "5.225"
$ws.Range("E39").Value = "  -1.73%  "

# --- Row 40 ---
Set-TextValue $ws.Range("D40") "0.2170"
$ws.Range("E40").Value = "  +0.76%  "

# --- Row 42 ---
Set-TextValue $ws.Range("D42") "1.219"
$ws.Range("E42").Value = "  -1.70%  "

# --- Row 43 ---
Set-TextValue $ws.Range("D43") "8.100"
$ws.Range("E43").Value = "  -1.44%  "

# --- Row 44 ---
Set-TextValue $ws.Range("D44") "14.33"
$ws.Range("E44").Value = "  +1.41%  "

# --- Row 46 ---
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D46") "3.875"
$ws.Range("E46").Value = "  +0.24%  "

# --- Row 47 ---
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.6128"
$ws.Range("E47").Value = "  -1.83%  "

# --- Row 48 ---
Set-TextValue $ws.Range("D48") "128.36"
$ws.Range("E48").Value = "  -3.25%  "

# --- Row 49 ---
Set-TextValue $ws.Range("D49") "2.057"
$ws.Range("E49").Value = "  -0.20%  "

# --- Row 50 ---
Set-TextValue $ws.Range("D50") "1.179"
$ws.Range("E50").Value = "  -1.90%  "

# --- Row 51 ---
Set-TextValue $ws.Range("D51") "0.07116"
$ws.Range("E51").Value = "  -5.17%  "

# --- Row 41 (volume only) ---
$ws.Range("E41").Value = "  +0.95%  "

# --- Row 45 (volume only) ---
$ws.Range("E45").Value = "  -0.15%  "
